$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Force text number format for D (Material) and J (Data de fim) columns ---
# so that purely-numeric-looking strings are stored as text, matching the
# original column typing (MATNR / GLTRP are text fields in this sheet).
$dRange = $ws.Range("D3:D23")
$jRange = $ws.Range("J3:J23")
$dRange.NumberFormat = "@"
$jRange.NumberFormat = "@"

# --- Update data rows 3-23 ---
$ws.Range("B3").Value = 3820016515
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = "71000843"
$ws.Range("F3").Value = "ZRSE"
$ws.Range("G3").Value = 23
$ws.Range("J3").Value = "20220408"

$ws.Range("B4").Value = 3820016515
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = "71000949"
$ws.Range("F4").Value = "ZRSE"
$ws.Range("G4").Value = 23
$ws.Range("J4").Value = "20220408"

$ws.Range("B5").Value = 3820016515
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = "71000481"
$ws.Range("F5").Value = "ZRPL"
$ws.Range("G5").Value = 500
$ws.Range("J5").Value = "20220408"

$ws.Range("B6").Value = 3823001120
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = "70011415"
$ws.Range("F6").Value = "ZCMQ"
$ws.Range("G6").Value = 400
$ws.Range("J6").Value = "20211217"

$ws.Range("B7").Value = 3823001120
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = "72000281"
$ws.Range("F7").Value = "ZCMM"
$ws.Range("G7").Value = 400
$ws.Range("J7").Value = "20211217"

$ws.Range("B8").Value = 3823001120
$ws.Range("C8").Value = 20
$ws.Range("D8").Value = "72000470"
$ws.Range("F8").Value = "ZCEM"
$ws.Range("G8").Value = 400
$ws.Range("J8").Value = "20211217"

$ws.Range("B9").Value = 3823001120
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = "70025396"
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = 384
$ws.Range("J9").Value = "20211217"

$ws.Range("B10").Value = 3823001120
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = "70026680"
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = 768
$ws.Range("J10").Value = "20211217"

$ws.Range("B11").Value = 3823001120
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = "70026679"
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = 768
$ws.Range("J11").Value = "20211217"

$ws.Range("B12").Value = 3823001123
$ws.Range("C12").Value = 190
$ws.Range("D12").Value = "71000923"
$ws.Range("F12").Value = "ZRSE"
$ws.Range("G12").Value = 23
$ws.Range("J12").Value = "20220408"

$ws.Range("B13").Value = 3823001123
$ws.Range("C13").Value = 190
$ws.Range("D13").Value = "71000094"
$ws.Range("F13").Value = "ZRSE"
$ws.Range("G13").Value = 23
$ws.Range("J13").Value = "20220408"

$ws.Range("B14").Value = 3823001123
$ws.Range("C14").Value = 190
$ws.Range("D14").Value = "71000430"
$ws.Range("F14").Value = "ZRPL"
$ws.Range("G14").Value = 1500
$ws.Range("J14").Value = "20220408"

$ws.Range("B15").Value = 3823001123
$ws.Range("C15").Value = 180
$ws.Range("D15").Value = "71000914"
$ws.Range("F15").Value = "ZRSE"
$ws.Range("G15").Value = 27
$ws.Range("J15").Value = "20220408"

$ws.Range("B16").Value = 3823001123
$ws.Range("C16").Value = 180
$ws.Range("D16").Value = "71000073"
$ws.Range("F16").Value = "ZRSE"
$ws.Range("G16").Value = 27
$ws.Range("J16").Value = "20220408"

$ws.Range("B17").Value = 3823001123
$ws.Range("C17").Value = 180
$ws.Range("D17").Value = "71000423"
$ws.Range("F17").Value = "ZRPL"
$ws.Range("G17").Value = 2000
$ws.Range("J17").Value = "20220408"

$ws.Range("B18").Value = 3823001123
$ws.Range("C18").Value = 170
$ws.Range("D18").Value = "71000808"
$ws.Range("F18").Value = "ZRSE"
$ws.Range("G18").Value = 25
$ws.Range("J18").Value = "20220408"

$ws.Range("B19").Value = 3823001123
$ws.Range("C19").Value = 170
$ws.Range("D19").Value = "71000098"
$ws.Range("F19").Value = "ZRSE"
$ws.Range("G19").Value = 25
$ws.Range("J19").Value = "20220408"

$ws.Range("B20").Value = 3823001123
$ws.Range("C20").Value = 170
$ws.Range("D20").Value = "71000544"
$ws.Range("F20").Value = "ZRPL"
$ws.Range("G20").Value = 2100
$ws.Range("J20").Value = "20220408"

$ws.Range("B21").Value = 3823001123
$ws.Range("C21").Value = 160
$ws.Range("D21").Value = "71000908"
$ws.Range("F21").Value = "ZRSE"
$ws.Range("G21").Value = 7
$ws.Range("J21").Value = "20220408"

$ws.Range("B22").Value = 3823001123
$ws.Range("C22").Value = 160
$ws.Range("D22").Value = "71000067"
$ws.Range("F22").Value = "ZRSE"
$ws.Range("G22").Value = 7
$ws.Range("J22").Value = "20220408"

$ws.Range("B23").Value = 3823001123
$ws.Range("C23").Value = 160
$ws.Range("D23").Value = "71000460"
$ws.Range("F23").Value = "ZRPL"
$ws.Range("G23").Value = 500
$ws.Range("J23").Value = "20220408"

# Restore normal style on D/J columns (keeps them as text without altering look)
$dRange.Style = "Normal"
$jRange.Style = "Normal"

# --- Rows 24-30 no longer have data; clear them back to the blank-row pattern ---
$ws.Range("A24:N30").ClearContents()
